$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("# Amino acids" shifts to E, etc.)
$ws.Columns("D").Insert()

# New column D: "Avg length" = # Amino acids / # Sequences
$ws.Range("D1").Value = "Avg length"
$ws.Range("D2:D15").Formula = "=E2/C2"
$ws.Range("D2:D15").NumberFormat = "0"

# Match column width of newly inserted column to its neighbours
$ws.Columns("D").ColumnWidth = 9.9

# Row heights: whole used range now uses the compact 13.8pt row height
$ws.Rows("1:15").RowHeight = 13.8

# Restore selection / active cell
$ws.Range("H11").Select()
